# Fruta / hortaliza, semanal
# Insert the latest week's 3 new price rows into the "Plátano" data table
# right before the existing row 490, shifting all subsequent rows down by 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows at position 490 (old row 490 and everything below
# move down to 493, 494, ...)
$ws.Rows("490:492").Insert()

# New row 490: Barraganete / Primera, fecha 2022-05-13 (serial 44694)
$ws.Range("A490").Value = 10
$ws.Range("B490").Value = "Vega Modelo de Temuco"
$ws.Range("C490").Value = "La Araucanía"
$ws.Range("D490").Value = 44694
$ws.Range("E490").Value = 9
$ws.Range("F490").Value = "Fruta"
$ws.Range("G490").Value = 100108
$ws.Range("H490").Value = "Tropicales y subtropicales"
$ws.Range("I490").Value = 100108006
$ws.Range("J490").Value = "Plátano"
$ws.Range("K490").Value = "Barraganete"
$ws.Range("L490").Value = "Primera"
$ws.Range("M490").Value = 35
$ws.Range("N490").Value = 25000
$ws.Range("O490").Value = 25000
$ws.Range("P490").Value = 25000
$ws.Range("Q490").Value = "$/caja 20 kilos"
$ws.Range("R490").Value = "Ecuador"
$ws.Range("S490").Value = 1250
$ws.Range("T490").Value = 20

# New row 491: Sin especificar / Maduro, fecha 2022-05-13 (serial 44694)
$ws.Range("A491").Value = 10
$ws.Range("B491").Value = "Vega Modelo de Temuco"
$ws.Range("C491").Value = "La Araucanía"
$ws.Range("D491").Value = 44694
$ws.Range("E491").Value = 9
$ws.Range("F491").Value = "Fruta"
$ws.Range("G491").Value = 100108
$ws.Range("H491").Value = "Tropicales y subtropicales"
$ws.Range("I491").Value = 100108006
$ws.Range("J491").Value = "Plátano"
$ws.Range("K491").Value = "Sin especificar"
$ws.Range("L491").Value = "Maduro"
$ws.Range("M491").Value = 125
$ws.Range("N491").Value = 12000
$ws.Range("O491").Value = 12000
$ws.Range("P491").Value = 12000
$ws.Range("Q491").Value = "$/caja 20 kilos"
$ws.Range("R491").Value = "Ecuador"
$ws.Range("S491").Value = 600
$ws.Range("T491").Value = 20

# New row 492: Sin especificar / Pintón, fecha 2022-05-13 (serial 44694)
$ws.Range("A492").Value = 10
$ws.Range("B492").Value = "Vega Modelo de Temuco"
$ws.Range("C492").Value = "La Araucanía"
$ws.Range("D492").Value = 44694
$ws.Range("E492").Value = 9
$ws.Range("F492").Value = "Fruta"
$ws.Range("G492").Value = 100108
$ws.Range("H492").Value = "Tropicales y subtropicales"
$ws.Range("I492").Value = 100108006
$ws.Range("J492").Value = "Plátano"
$ws.Range("K492").Value = "Sin especificar"
$ws.Range("L492").Value = "Pintón"
$ws.Range("M492").Value = 930
$ws.Range("N492").Value = 14000
$ws.Range("O492").Value = 15000
$ws.Range("P492").Value = 14269
$ws.Range("Q492").Value = "$/caja 20 kilos"
$ws.Range("R492").Value = "Ecuador"
$ws.Range("S492").Value = 713
$ws.Range("T492").Value = 20
